$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task rows: Refresh Token / AutoMapper packages and Confirm Email task
$ws.Range("B11").Value = "Microsoft.AspNetCore.Mvc.NewtonsoftJson"
$ws.Range("B12").Value = "AutoMapper, AutoMapper.Extensions.Microsoft.DependencyInjection"
$ws.Range("A13").Value = "Sign in Email confirmation"
$ws.Range("B13").Value = "Add .SignIn.RequireConfirmedEmail = true; to services.AddIdentity,`nInstall NETCore.MailKit, Configure MailKitOptions at Startup.cs file, Add Email setting in the appSettings.json file."

# B13 wraps like the other "instructions" column cells
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = -4160
$ws.Rows(13).RowHeight = 43.2

# Column B got a bit wider to fit the new text
$ws.Columns("B").ColumnWidth = 56.5

# Update view state: scrolled down a couple rows, selection now on B13
$ws.Range("B13").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
